# 3.1 Example 2 output from YAML now matches that from CTI.
# Rename the sheet so it reflects the "Example 2" dataset.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "ceO2_output_rates_gas_Example2"

# Keep the embedded scatter chart's series title in sync with the renamed
# sheet (the chart lives on the same sheet and referenced the old name).
$co = $ws.ChartObjects(1)
$ser = $co.Chart.SeriesCollection(1)
$ser.Formula = "=SERIES(ceO2_output_rates_gas_Example2!`$C`$1,ceO2_output_rates_gas_Example2!`$B`$2:`$B`$501,ceO2_output_rates_gas_Example2!`$C`$2:`$C`$501,1)"

# Match the author's last selected cell on the sheet.
$ws.Range("F29").Select()
